$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.243.82"
$ws.Range("E2").Value = "  -3.20%  "

$ws.Range("D3").Value = "3.696.74"
$ws.Range("E3").Value = "  -1.89%  "

$ws.Range("E4").Value = "  +0.06%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "592.34"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -3.64%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "166.32"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.97%  "

$ws.Range("D7").Value = "3.693.84"
$ws.Range("E7").Value = "  -1.98%  "

$ws.Range("E8").Value = "  +0.02%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.522"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.72%  "

$ws.Range("E10").Value = "  -2.57%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.14"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.93%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.459"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -5.14%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "37.70"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -5.06%  "

$ws.Range("E14").Value = "  -5.10%  "

$ws.Range("D15").Value = "4.318.86"
$ws.Range("E15").Value = "  -1.70%  "

$ws.Range("D16").Value = "3.696.73"
$ws.Range("E16").Value = "  -1.74%  "

$ws.Range("D17").Value = "67.261.20"
$ws.Range("E17").Value = "  -3.27%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.114"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -3.80%  "

$ws.Range("E19").Value = "  -6.07%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.14"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.58%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "487.55"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.30%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "9.12"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.56%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.720"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.76%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "84.96"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.59%  "

$ws.Range("E25").Value = "  -6.41%  "

$ws.Range("E26").Value = "  +0.13%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "12.13"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -5.93%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.95"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -5.51%  "

$ws.Range("E30").Value = "  -2.78%  "

$ws.Range("E31").Value = "  -6.11%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.70"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -5.63%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "31.60"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.70%  "

$ws.Range("E34").Value = "  -7.27%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.02%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -4.77%  "

$ws.Range("E37").Value = "  -6.49%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.131"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -7.39%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.321"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -5.54%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "444.51"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -6.47%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "48.83"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.86%  "

$ws.Range("E42").Value = "  -5.43%  "

$ws.Range("E43").Value = "  -7.15%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.27"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.66%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "39.84"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -10.13%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "140.58"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.73%  "

$ws.Range("D48").Value = "2.779.35"
$ws.Range("E48").Value = "  -5.52%  "

$ws.Range("E49").Value = "  -4.68%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "25.26"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -8.37%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "23.60"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +7.34%  "
